# Update the FilesTab Neo4j/Cypher query in cell B4 of the "startup" sheet.
# The edit removes the `File Type` column (f.file_type) and the `Breed`
# column (demo.breed) from the RETURN clause, matching the upstream
# "10 icdc scripts for jenkins" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = "MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['GLIOMA01'] 
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS ``File Name``, 
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$cell = $ws.Range("B4")
$cell.Value = $newQuery

# The author's selection moved from C4 to B4 after editing the cell.
$ws.Activate()
$cell.Select()
